# "Correction : Erreur Changement de date"
#
# The sprint calendar's header row (row 2) was driven by volatile
# TODAY()+/-n formulas, which made the two-week date strip drift every
# time the workbook was opened/recalculated. This pins the strip to a
# fixed start date (2024-02-05) and updates the "task completed on ..."
# labels (rows 21-27) to the correct mid-February dates instead of the
# late-February/March ones that had drifted in.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# --- Row 2: replace the volatile TODAY()-based header dates with a
#     fixed two-week strip starting Monday 2024-02-05 ---
$ws.Range("B2").Formula = "=DATE(2024,2,5)"
$ws.Range("C2").Formula = "=B2+1"
$ws.Range("D2:O2").Formula = "=C2+1"

# --- Rows 21-27: correct the "task completed" legend date labels ---
$ws.Range("B21").Value = "Tâches terminé le 06/02/2024"
$ws.Range("B22").Value = "Tâche terminé le 09/02/2024"
$ws.Range("B23").Value = "Tâches terminé le 12/02/2024"
$ws.Range("B24").Value = "Tâches terminé le 13/02/2024"
$ws.Range("B25").Value = "Tâche terminé le 14/02/2024"
$ws.Range("B26").Value = "Tâche terminé le 16/02/2024"
$ws.Range("B27").Value = "Tâches terminé le 18/02/2024"

# --- Drop the stray formatted-but-empty cell at E25 ---
$ws.Range("E25").Clear()

# --- Add the new empty, date-formatted cell A30 (same number format as
#     the row-2 date strip) ---
$ws.Range("B2").Copy()
$ws.Range("A30").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Restore the active selection to D10 ---
$ws.Range("D10").Select()
